# Established connection to pumps
# Only Pump 1 (row 2) stays configured/connected; the previously configured
# pumps 2-5 (rows 3-6, columns A-H) are reset back to blank (values cleared,
# formatting/styles kept), and the now-unused pump id labels disappear from
# the shared string table as a natural consequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the pump id / pH settings for rows 3-6 (pumps 2-5), keeping the
# existing cell formatting untouched.
$ws.Range("A3:H6").ClearContents()

# Move/save the current selection, matching where the user left off.
$ws.Range("U10").Select()
